$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.43 = 30153.46 pesos`n✅ 30153.46 pesos = 7.39 = 960.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 134.5
$wsTasas.Range("O10").Value = 4055.64
$wsTasas.Range("N12").Value = 4082
$wsTasas.Range("O12").Value = 130.001
